$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.962.97"
$ws.Range("E2").Value = "  +2.01%  "

Set-TextValue $ws.Range("D3") "1.904.72"
$ws.Range("E3").Value = "  +1.94%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "333.01"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("E6").Value = "  -0.02%  "

Set-TextValue $ws.Range("D7") "0.4642"
$ws.Range("E7").Value = "  -1.22%  "

Set-TextValue $ws.Range("D8") "0.4058"
$ws.Range("E8").Value = "  +2.54%  "

Set-TextValue $ws.Range("D9") "47.95"
$ws.Range("E9").Value = "  +1.27%  "

Set-TextValue $ws.Range("D10") "0.08000"
$ws.Range("E10").Value = "  -0.10%  "

Set-TextValue $ws.Range("D11") "1.001"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("E12").Value = "  -0.88%  "

Set-TextValue $ws.Range("D13") "1.910.21"
$ws.Range("E13").Value = "  +2.58%  "

Set-TextValue $ws.Range("D14") "5.915"
$ws.Range("E14").Value = "  -1.20%  "

Set-TextValue $ws.Range("D15") "7.058"
$ws.Range("E15").Value = "  -2.31%  "

$ws.Range("E16").Value = "  -0.11%  "

Set-TextValue $ws.Range("D17") "88.89"
$ws.Range("E17").Value = "  -2.61%  "

$ws.Range("E18").Value = "  -0.75%  "

Set-TextValue $ws.Range("D19") "0.06558"
$ws.Range("E19").Value = "  -1.11%  "

Set-TextValue $ws.Range("D20") "17.39"
$ws.Range("E20").Value = "  -0.94%  "

Set-TextValue $ws.Range("D21") "1.002"
$ws.Range("E21").Value = "  +0.28%  "

Set-TextValue $ws.Range("D22") "28.973.84"
$ws.Range("E22").Value = "  +1.98%  "

Set-TextValue $ws.Range("D23") "5.454"
$ws.Range("E23").Value = "  +0.04%  "

Set-TextValue $ws.Range("D24") "11.17"
$ws.Range("E24").Value = "  +1.29%  "

Set-TextValue $ws.Range("D25") "2.238"
$ws.Range("E25").Value = "  -1.34%  "

Set-TextValue $ws.Range("D26") "2.130.32"
$ws.Range("E26").Value = "  +2.14%  "

Set-TextValue $ws.Range("D27") "157.66"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("E28").Value = "  -0.34%  "

Set-TextValue $ws.Range("D29") "2.095"
$ws.Range("E29").Value = "  -1.56%  "

Set-TextValue $ws.Range("D30") "5.390"
$ws.Range("E30").Value = "  -1.98%  "

Set-TextValue $ws.Range("D31") "118.75"
$ws.Range("E31").Value = "  -1.06%  "

Set-TextValue $ws.Range("D32") "0.9800"
$ws.Range("E32").Value = "  +1.54%  "

Set-TextValue $ws.Range("D33") "0.09374"
$ws.Range("E33").Value = "  -1.12%  "

Set-TextValue $ws.Range("D34") "1.414"
$ws.Range("E34").Value = "  +3.02%  "

Set-TextValue $ws.Range("D35") "3.603"
$ws.Range("E35").Value = "  +0.91%  "

Set-TextValue $ws.Range("D36") "5.287"
$ws.Range("E36").Value = "  -1.13%  "

Set-TextValue $ws.Range("D37") "0.06073"
$ws.Range("E37").Value = "  -0.16%  "

Set-TextValue $ws.Range("D38") "0.02222"
$ws.Range("E38").Value = "  -0.91%  "

Set-TextValue $ws.Range("D39") "8.395"
$ws.Range("E39").Value = "  +0.28%  "

Set-TextValue $ws.Range("D40") "1.163"
$ws.Range("E40").Value = "  -1.92%  "

Set-TextValue $ws.Range("D41") "1.002"
$ws.Range("E41").Value = "  +0.15%  "

Set-TextValue $ws.Range("D42") "0.5787"
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("E43").Value = "  -2.47%  "

Set-TextValue $ws.Range("D44") "10.10"
$ws.Range("E44").Value = "  -2.14%  "

Set-TextValue $ws.Range("D45") "1.266"
$ws.Range("E45").Value = "  -1.94%  "

Set-TextValue $ws.Range("D46") "2.322"
$ws.Range("E46").Value = "  +12.49%  "

Set-TextValue $ws.Range("D47") "12.09"
$ws.Range("E47").Value = "  -0.71%  "

Set-TextValue $ws.Range("D48") "0.5487"
$ws.Range("E48").Value = "  -1.59%  "

Set-TextValue $ws.Range("D49") "1.899"
$ws.Range("E49").Value = "  -2.75%  "

Set-TextValue $ws.Range("D50") "0.07016"
$ws.Range("E50").Value = "  +2.34%  "

Set-TextValue $ws.Range("D51") "47.34"
$ws.Range("E51").Value = "  +22.44%  "
